$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 into H1 so the new "Save" header picks up
# the same bold/bordered header formatting (style index 1), then set text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New "Save" data column (0/1 flags) for rows 2-31.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
